$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Shartnoma raqam), G (Telefon raqam) and H (Sana) must be stored as
# literal text (matching the source data which is all inline-string), not as
# numbers/dates, even though their contents look numeric / date-like.
$ws.Range("D157:D172").NumberFormat = "@"
$ws.Range("G157:G172").NumberFormat = "@"
$ws.Range("H157:H172").NumberFormat = "@"

# Row 157
$ws.Range("A157").Value = 'Aleksevnina Elvira Samigullayevna'
$ws.Range("B157").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C157").Value = 'AD1153885'
$ws.Range("D157").Value = '332'
$ws.Range("E157").Value = 'Toshkent shahri'
$ws.Range("F157").Value = 'Yunusobod tumani'
$ws.Range("G157").Value = '998903296045'
$ws.Range("H157").Value = '03-12-2024'

# Row 158
$ws.Range("A158").Value = 'Karimova Nigora Abdug''oppor qizi'
$ws.Range("B158").Value = 'Maktabgacha talim tashkiloti defektolog/logopedi'
$ws.Range("C158").Value = 'AB1487410'
$ws.Range("D158").Value = '333'
$ws.Range("E158").Value = 'Fargona viloyati'
$ws.Range("F158").Value = 'Fargʻona tumani'
$ws.Range("G158").Value = '998911268606'
$ws.Range("H158").Value = '03-12-2024'

# Row 159
$ws.Range("A159").Value = 'Nabiyeva Diyora Mahmud qizi'
$ws.Range("B159").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C159").Value = 'AA9586555'
$ws.Range("D159").Value = '334'
$ws.Range("E159").Value = 'Toshkent shahri'
$ws.Range("F159").Value = 'Yunusobod tumani'
$ws.Range("G159").Value = '998773434046'
$ws.Range("H159").Value = '03-12-2024'

# Row 160
$ws.Range("A160").Value = 'Yuldasheva ShaxodatAbdugani Qizi'
$ws.Range("B160").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C160").Value = 'AB3849334'
$ws.Range("D160").Value = '335'
$ws.Range("E160").Value = 'Namangan viloyati'
$ws.Range("F160").Value = 'Yangi Namangan'
$ws.Range("G160").Value = '998940084991'
$ws.Range("H160").Value = '04-12-2024'

# Row 161
$ws.Range("A161").Value = 'Xolmatova Nazokat Abdurashidovna'
$ws.Range("B161").Value = 'Maktabgacha talim tashkiloti tarbiyachisi 576 soat'
$ws.Range("C161").Value = 'AD3224685'
$ws.Range("D161").Value = '336'
$ws.Range("E161").Value = 'Toshkent viloyati'
$ws.Range("F161").Value = 'Yangiyoʻl tumani'
$ws.Range("G161").Value = '998944176679'
$ws.Range("H161").Value = '04-12-2024'

# Row 162
$ws.Range("A162").Value = 'Qarshiyeva Nargiza'
$ws.Range("B162").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C162").Value = 'AA6933915'
$ws.Range("D162").Value = '337'
$ws.Range("E162").Value = 'Qashqadaryo viloyati'
$ws.Range("F162").Value = 'Yakkabogʻ tumani'
$ws.Range("G162").Value = '998933948706'
$ws.Range("H162").Value = '04-12-2024'

# Row 163
$ws.Range("A163").Value = 'Jalolova Shahnoza  Nizomiddinovna'
$ws.Range("B163").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C163").Value = 'AB6670685'
$ws.Range("D163").Value = '338'
$ws.Range("E163").Value = 'Toshkent shahri'
$ws.Range("F163").Value = 'Chilonzor tumani'
$ws.Range("G163").Value = '998917909191'
$ws.Range("H163").Value = '04-12-2024'

# Row 164
$ws.Range("A164").Value = 'Normatova Zamira Latifovna'
$ws.Range("B164").Value = 'Maktabgacha talim tashkiloti direktori'
$ws.Range("C164").Value = 'AD2136078'
$ws.Range("D164").Value = '339'
$ws.Range("E164").Value = 'Andijon viloyati'
$ws.Range("F164").Value = 'Asaka tumani'
$ws.Range("G164").Value = '998934122687'
$ws.Range("H164").Value = '05-12-2024'

# Row 165
$ws.Range("A165").Value = 'G''oyibova Dilshoda Hamroqulovna'
$ws.Range("B165").Value = 'Maktabgacha talim tashkiloti metodisti'
$ws.Range("C165").Value = 'AD8546284'
$ws.Range("D165").Value = '340'
$ws.Range("E165").Value = 'Buxoro viloyati'
$ws.Range("F165").Value = 'Qorakoʻl tumani'
$ws.Range("G165").Value = '998930820105'
$ws.Range("H165").Value = '05-12-2024'

# Row 166
$ws.Range("A166").Value = 'Raxmonova Dilafruz Ochilovna'
$ws.Range("B166").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C166").Value = 'AB3464045'
$ws.Range("D166").Value = '341'
$ws.Range("E166").Value = 'Buxoro viloyati'
$ws.Range("F166").Value = 'Qorakoʻl tumani'
$ws.Range("G166").Value = '998943280180'
$ws.Range("H166").Value = '05-12-2024'

# Row 167
$ws.Range("A167").Value = 'Tadjiyeva Svetlana Aktyamovna'
$ws.Range("B167").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C167").Value = 'AD2327516'
$ws.Range("D167").Value = '342'
$ws.Range("E167").Value = 'Toshkent shahri'
$ws.Range("F167").Value = 'Yunusobod tumani'
$ws.Range("G167").Value = '998948364949'
$ws.Range("H167").Value = '05-12-2024'

# Row 168
$ws.Range("A168").Value = 'Normatova Zamiraxon Latifovna'
$ws.Range("B168").Value = 'Maktabgacha talim tashkiloti tarbiyachisi 576 soat'
$ws.Range("C168").Value = 'AD2136078'
$ws.Range("D168").Value = '343'
$ws.Range("E168").Value = 'Andijon viloyati'
$ws.Range("F168").Value = 'Asaka tumani'
$ws.Range("G168").Value = '998934122687'
$ws.Range("H168").Value = '05-12-2024'

# Row 169
$ws.Range("A169").Value = 'Mutalova Nasiba Ochilova'
$ws.Range("B169").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C169").Value = 'AB9594554'
$ws.Range("D169").Value = '344'
$ws.Range("E169").Value = 'Buxoro viloyati'
$ws.Range("F169").Value = 'Qorakoʻl tumani'
$ws.Range("G169").Value = '998939689434'
$ws.Range("H169").Value = '05-12-2024'

# Row 170
$ws.Range("A170").Value = 'Saliea GulnozaTursunbayvna'
$ws.Range("B170").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C170").Value = 'AB4918887'
$ws.Range("D170").Value = '345'
$ws.Range("E170").Value = 'Toshkent shahri'
$ws.Range("F170").Value = 'Olmazor tumani'
$ws.Range("G170").Value = '998916479256'
$ws.Range("H170").Value = '06-12-2024'

# Row 171
$ws.Range("A171").Value = 'Saliea GulnozaTursunbayvna'
$ws.Range("B171").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C171").Value = 'AB4918887'
$ws.Range("D171").Value = '346'
$ws.Range("E171").Value = 'Toshkent shahri'
$ws.Range("F171").Value = 'Olmazor tumani'
$ws.Range("G171").Value = '998916479256'
$ws.Range("H171").Value = '06-12-2024'

# Row 172
$ws.Range("A172").Value = 'Fofurova Dilrabo Ravshanbekovna'
$ws.Range("B172").Value = 'Maktabgacha talim tashkiloti tarbiyachisi'
$ws.Range("C172").Value = 'AD1844351'
$ws.Range("D172").Value = '347'
$ws.Range("E172").Value = 'Toshkent shahri'
$ws.Range("F172").Value = 'Yunusobod tumani'
$ws.Range("G172").Value = '998998859987'
$ws.Range("H172").Value = '06-12-2024'

